$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unraveling the Enigmatic Beauty: Quantum Entanglement",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "The Enigmatic World of Chemistry: Unveiling the Secrets of Matter", $wdReplaceAll)

# ---------------------------------------------------------------------
# 2. Author name paragraph: "Dr. Chloe Anderson" (3 runs) -> "Olivia Curtis"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dr. Chloe Anderson",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Olivia Curtis", $wdReplaceAll)

# ---------------------------------------------------------------------
# 3. Email paragraph
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "chloe",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "oliviacurtis", $wdReplaceAll)

$d.Content.Find.Execute(
    "anderson@quantumresearch",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "chemistry@validmail", $wdReplaceAll)

# ---------------------------------------------------------------------
# 4. First body paragraph
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "In the realm of quantum physics, where uncertainty reigns and particles defy classical logic, lies a captivating phenomenon known as quantum entanglement",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "In the vast expanse of the universe, chemistry stands as a pillar of understanding, unveiling the fundamental principles that govern the intricate interactions of matter at its most basic level",
    $wdReplaceAll)

# " This extraordinary phenomenon..." sentence is replaced, and two more
# sentences are appended right after it (new runs).
$rng = $d.Content
$rng.Find.Execute(
    "This extraordinary phenomenon, often labeled as " + [char]34 + "spooky action at a distance," + [char]34 + " by Albert Einstein, challenges our intuitive understanding of reality and opens up a new realm of scientific exploration")
$rng.Text = "It embarks upon an exploration of the diverse elements and compounds that make up our world, deciphering the enigmatic language of chemical reactions and revealing the profound impact they have on our existence"
$rng.Collapse(0)
$rng.InsertAfter(". Chemistry weaves its way through fabrics of our lives, touching every aspect from our clothes, and medicine to our food and technology, unlocking the secrets of the molecular dance that shapes our reality")

# "Quantum entanglement is a captivating dance..." -> "Introduction Continued:"
# followed by a blank line and a new paragraph of body text.
$rng = $d.Content
$rng.Find.Execute("Quantum entanglement is a captivating dance between two or more particles whose properties, such as spin, polarization, or energy, become correlated in an inexplicable manner")
$rng.Text = "Introduction Continued:"
$rng.Collapse(0)
$rng.InsertAfter([char]11)
$rng.Collapse(0)
$rng.InsertAfter([char]11)
$rng.Collapse(0)
$rng.InsertAfter("The history of chemistry is a rich tapestry woven with tales of intrepid pioneers who dared to question the enigmatic nature of matter")

$d.Content.Find.Execute(
    "These particles remain intimately connected, regardless of the distance separating them, sharing their fates in a profound and nonlocal way",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "From ancient alchemists seeking the elusive philosopher's stone to modern chemists unraveling the intricate secrets of DNA, the quest for knowledge has driven the evolution of this field",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "The measurement of the property of one entangled particle instantaneously influences the state of the other, irrespective of the vast cosmic expanse separating them",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Each discovery, each breakthrough, has added a brushstroke of color to the vibrant canvas of our understanding, revealing the hidden harmonies of the chemical world",
    $wdReplaceAll)

# "Einstein famously declared..." -> "Introduction Continued:" + blank line + new paragraph
$rng = $d.Content
$rng.Find.Execute("Einstein famously declared this phenomenon as " + [char]34 + "spooky action at a distance," + [char]34 + " highlighting its ethereal and counterintuitive nature")
$rng.Text = "Introduction Continued:"
$rng.Collapse(0)
$rng.InsertAfter([char]11)
$rng.Collapse(0)
$rng.InsertAfter([char]11)
$rng.Collapse(0)
$rng.InsertAfter("Chemistry is not merely a body of knowledge; it is an art form, a symphony of elements and compounds that harmonize in countless ways")

# "Quantum entanglement has ignited..." -> several new sentences
$rng = $d.Content
$rng.Find.Execute("Quantum entanglement has ignited intense scientific debates, challenging our fundamental understanding of physics and pushing the boundaries of human knowledge")
$rng.Text = "It is the molecular dance of atoms, the delicate interplay of forces that shape the world around us"
$rng.Collapse(0)
$rng.InsertAfter(". Chemistry teaches us to observe, to experiment, to analyze, to create")
$rng.Collapse(0)
$rng.InsertAfter(". It cultivates a sense of wonder and curiosity, encouraging us to ask questions about the world around us and to seek answers in the intricate web of chemical processes that underlie all of nature")

# ---------------------------------------------------------------------
# 5. Summary paragraph
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Quantum entanglement stands as an enigmatic masterpiece of nature, a symphony of subatomic choreography that weaves an intricate tapestry of interconnectedness",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Chemistry delves into the enigmatic secrets of matter, unveiling the profound impact that chemical reactions have on our lives",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "It redefines our notions of locality and causality, ushering in a profound transformation in how we perceive the universe",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "It is a tapestry woven with the threads of history, where the quest for knowledge has driven the evolution of this field",
    $wdReplaceAll)

# Remove the "While its underlying..." / "." / "From cryptography..." runs
# and replace with a single new sentence.
$rng = $d.Content
$rng.Find.Execute("While its underlying mechanisms may still elude our grasp, quantum entanglement holds the promise of unlocking transformative technologies and reshaping our understanding of the quantum realm. From cryptography to computation, the implications of quantum entanglement reverberate across diverse fields, beckoning us to explore this extraordinary phenomenon and its captivating implications for science, technology, and human knowledge")
$rng.Text = "Chemistry is not merely a body of knowledge, but an art form, a symphony of elements and compounds that harmonize in countless ways, revealing the hidden harmonies of the chemical world"

# ---------------------------------------------------------------------
# 6. New trailing empty paragraph before the section break.
# ---------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertBefore([char]13)
